# Atualização de bases das ligas, do dia: 17-05-2024 às 13:59
# Swap full row contents (columns B..AB) between pairs of adjacent rows
# whose "id" values (column B) got reordered in the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($row1, $row2, $colStart, $colEnd) {
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $cell1 = $ws.Cells.Item($row1, $c)
        $cell2 = $ws.Cells.Item($row2, $c)
        $tmp = $cell1.Value2
        $cell1.Value = $cell2.Value2
        $cell2.Value = $tmp
    }
}

# Row pairs (1-based sheet rows) whose B:AB data must be swapped.
# Column A (row index) is left untouched for every pair.
$pairs = @(
    @(38, 39),
    @(49, 50),
    @(65, 66),
    @(71, 72),
    @(115, 116),
    @(125, 126),
    @(128, 129),
    @(146, 147),
    @(191, 192),
    @(248, 249)
)

foreach ($p in $pairs) {
    Swap-RowData $p[0] $p[1] 2 28
}
